$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Add a new worksheet positioned after Sheet1
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

# Populate Sheet2 with data
$ws2.Range("A1").Value = "Another"
$ws2.Range("B1").Value = "Sheet"
$ws2.Range("A2").Value = "A"
$ws2.Range("B2").Value = "S"

# Update selection on the new sheet to match the target (A3)
$ws2.Range("A3").Select()
